$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 76, pushing the existing rows (old 76..100)
# down to become rows 78..102.
$ws.Rows.Item(76).Insert()
$ws.Rows.Item(76).Insert()

# Row 76: new weekly price entry (copy the constant columns from the row
# just below, which still holds the former row-76 data at this point).
$ws.Cells.Item(76, 1).Value = $ws.Cells.Item(78, 1).Value2
$ws.Cells.Item(76, 2).Value = $ws.Cells.Item(78, 2).Value2
$ws.Cells.Item(76, 3).Value = $ws.Cells.Item(78, 3).Value2
$ws.Cells.Item(76, 4).Value = 44736
$ws.Cells.Item(76, 5).Value = $ws.Cells.Item(78, 5).Value2
$ws.Cells.Item(76, 6).Value = $ws.Cells.Item(78, 6).Value2
$ws.Cells.Item(76, 7).Value = $ws.Cells.Item(78, 7).Value2
$ws.Cells.Item(76, 8).Value = $ws.Cells.Item(78, 8).Value2
$ws.Cells.Item(76, 9).Value = $ws.Cells.Item(78, 9).Value2
$ws.Cells.Item(76, 10).Value = $ws.Cells.Item(78, 10).Value2
$ws.Cells.Item(76, 11).Value = "Hass"
$ws.Cells.Item(76, 12).Value = "Primera"
$ws.Cells.Item(76, 13).Value = 400
$ws.Cells.Item(76, 14).Value = 16000
$ws.Cells.Item(76, 15).Value = 17000
$ws.Cells.Item(76, 16).Value = 16500
$ws.Cells.Item(76, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(76, 18).Value = "Perú"
$ws.Cells.Item(76, 19).Value = 1650
$ws.Cells.Item(76, 20).Value = 10

# Row 77: second new weekly price entry.
$ws.Cells.Item(77, 1).Value = $ws.Cells.Item(78, 1).Value2
$ws.Cells.Item(77, 2).Value = $ws.Cells.Item(78, 2).Value2
$ws.Cells.Item(77, 3).Value = $ws.Cells.Item(78, 3).Value2
$ws.Cells.Item(77, 4).Value = 44736
$ws.Cells.Item(77, 5).Value = $ws.Cells.Item(78, 5).Value2
$ws.Cells.Item(77, 6).Value = $ws.Cells.Item(78, 6).Value2
$ws.Cells.Item(77, 7).Value = $ws.Cells.Item(78, 7).Value2
$ws.Cells.Item(77, 8).Value = $ws.Cells.Item(78, 8).Value2
$ws.Cells.Item(77, 9).Value = $ws.Cells.Item(78, 9).Value2
$ws.Cells.Item(77, 10).Value = $ws.Cells.Item(78, 10).Value2
$ws.Cells.Item(77, 11).Value = "Hass"
$ws.Cells.Item(77, 12).Value = "Segunda"
$ws.Cells.Item(77, 13).Value = 400
$ws.Cells.Item(77, 14).Value = 15000
$ws.Cells.Item(77, 15).Value = 16000
$ws.Cells.Item(77, 16).Value = 15500
$ws.Cells.Item(77, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(77, 18).Value = "Perú"
$ws.Cells.Item(77, 19).Value = 1550
$ws.Cells.Item(77, 20).Value = 10
